$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '36.260.59'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.044.05'
$ws.Range('E3').Value = '  -2.63%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '244.67'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').Value = '56.56'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = '63.26'
$ws.Range('E9').Value = '  +5.74%  '
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  -3.49%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('D13').Value = '0.909'
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').Value = '14.13'
$ws.Range('E14').Value = '  -6.49%  '
$ws.Range('D15').Value = '2.335.74'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('E16').Value = '  -3.96%  '
$ws.Range('D17').Value = '2.020.53'
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').Value = '17.50'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '36.224.28'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').Value = '71.27'
$ws.Range('E20').Value = '  -2.99%  '
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('D22').Value = '236.98'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('E23').Value = '  -6.10%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('D27').Value = '9.28'
$ws.Range('E27').Value = '  -6.49%  '
$ws.Range('D28').Value = '164.38'
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').Value = '19.92'
$ws.Range('E29').Value = '  -5.04%  '
$ws.Range('E30').Value = '  -2.52%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  -8.26%  '
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = '4.39'
$ws.Range('E34').Value = '  -7.57%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.0868'
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').Value = '2.21'
$ws.Range('E38').Value = '  -9.46%  '
$ws.Range('D39').Value = '5.05'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('D40').Value = '1.21'
$ws.Range('E40').Value = '  -5.80%  '
$ws.Range('D41').Value = '2.88'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('D42').Value = '0.0215'
$ws.Range('E42').Value = '  -3.26%  '
$ws.Range('E43').Value = '  -6.69%  '
$ws.Range('D44').Value = '93.19'
$ws.Range('E44').Value = '  -3.94%  '
$ws.Range('D45').Value = '0.0905'
$ws.Range('E45').Value = '  -5.81%  '
$ws.Range('D46').Value = '15.89'
$ws.Range('E46').Value = '  -3.53%  '
$ws.Range('D47').Value = '1.365.78'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('D48').Value = '7.37'
$ws.Range('E48').Value = '  +3.69%  '
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').Value = '  -6.80%  '
$ws.Range('D51').Value = '45.66'
$ws.Range('E51').Value = '  -0.80%  '
